$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29: update % Build to 0.67 ---
$ws.Range("D29").Value = 0.67

# --- Row 30: add Finish date, mark % Build complete ---
# Match the date-cell formatting already used in column B30 (style s="25")
$ws.Range("B30").Copy()
$ws.Range("C30").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C30").Value = 44474
$ws.Range("D30").Value = 1

# --- Row 31: add Start/Finish dates, mark % Build complete ---
$ws.Range("B30").Copy()
$ws.Range("B31").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B31").Value = 44474

$ws.Range("B30").Copy()
$ws.Range("C31").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C31").Value = 44474

$ws.Range("D31").Value = 1

$excel.CutCopyMode = 0

# --- Update the active sheet view / selection to match the edit location ---
$ws.Range("E32").Select()
$excel.ActiveWindow.ScrollRow = 25
